$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edits 1-3: merge several split runs back into single runs inside the
# paragraph "... cu 5: 3 15 10 (if n mod 3 ˄(AND) n mod 5)" — i.e. collapse
# " " + "1" + "5 " -> " 15 ", " " + "(if n mod 3 " -> " (if n mod 3 ", and
# "(" + "AND" + ") n mod 5)" -> "(AND) n mod 5)", while leaving the
# surrounding "10" (struck-through) and "˄" runs untouched / unmerged.
# ---------------------------------------------------------------------------

# --- Edit 1: " " + "1" + "5 " -> " 15 " --------------------------------
# Scope the search to start after "cu 3 " so it can't match the similar
# "... cu 3 sau cu 5: 3 5 6 9 11 ..." paragraph above it.
$searchRange1 = $d.Range(440, $d.Content.End)
$searchRange1.Find.Execute(" 15 ", $true, $false, $false, $false, $false, `
    $true, 1, $false, " 15 ", 2) | Out-Null

# --- Edit 2: " " + "(if n mod 3 " -> " (if n mod 3 " ---------------------
# Scope the search to start right after the (struck-through) "10" run so it
# can't match the other "... 11 (if n mod 3 ˅(OR) n mod 5)" paragraph.
$searchRange2 = $d.Range(495, $d.Content.End)
$searchRange2.Find.Execute(" (if n mod 3 ", $true, $false, $false, $false, `
    $false, $true, 1, $false, " (if n mod 3 ", 2) | Out-Null

# --- Edit 3: "(" + "AND" + ") n mod 5)" -> "(AND) n mod 5)" --------------
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("˄(AND) n mod 5)", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$symStart = $anchor.Start
$parenStart = $anchor.Start + 1
$searchRange3 = $d.Range($parenStart, $d.Content.End)
$searchRange3.Find.Execute("(AND) n mod 5)", $true, $false, $false, $false, `
    $false, $true, 1, $false, "(AND) n mod 5)", 2) | Out-Null

# The replace above leaves "˄" and "(AND) n mod 5)" coalesced into a single
# run (since they share identical formatting); re-split them apart by
# nudging a character property on the "˄" position only.
$symRange = $d.Range($symStart, $parenStart)
$symRange.Font.Bold = 1
$symRange.Font.Bold = 0

# ---------------------------------------------------------------------------
# Edit 4: split the run " p ˄ q" into " p ˄ " + "˥" + "q" (same formatting).
# ---------------------------------------------------------------------------
$anchor4 = $d.Content.Duplicate
$anchor4.Find.Execute(" p ˄ q", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$qPos = $anchor4.End - 1

# Insert the new "˥" character right before the final "q".
$ins = $d.Range($qPos, $qPos)
$ins.InsertAfter("˥")

# Force a run boundary between "˥" and "q" by nudging (and restoring) a
# character property on just the "q" run.
$qRange = $d.Range($qPos + 1, $qPos + 2)
$qRange.Font.Bold = 1
$qRange.Font.Bold = 0

# Force a run boundary between " p ˄ " and "˥" the same way.
$symRange4 = $d.Range($qPos, $qPos + 1)
$symRange4.Font.Bold = 1
$symRange4.Font.Bold = 0

Write-Output "Done"
